$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 125, pushing existing rows 125.. down by one.
$ws.Rows.Item(125).Insert()

# Populate the newly inserted row 125 with the new weekly record.
$ws.Cells.Item(125, 1).Value = 3
$ws.Cells.Item(125, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(125, 3).Value = "Coquimbo"
$ws.Cells.Item(125, 4).Value = 44754
$ws.Cells.Item(125, 5).Value = 5
$ws.Cells.Item(125, 6).Value = 100112026
$ws.Cells.Item(125, 7).Value = "Haba"
$ws.Cells.Item(125, 8).Value = "Sin especificar"
$ws.Cells.Item(125, 9).Value = "Primera"
$ws.Cells.Item(125, 10).Value = 85
$ws.Cells.Item(125, 11).Value = 17000
$ws.Cells.Item(125, 12).Value = 18000
$ws.Cells.Item(125, 13).Value = 17529
$ws.Cells.Item(125, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(125, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(125, 16).Value = 701
$ws.Cells.Item(125, 17).Value = 25
$ws.Cells.Item(125, 18).Value = "Hortaliza"
